$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $cell = $ws.Range($cellAddr)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue 'D2' '29.814.03'
Set-TextValue 'E2' '  -1.48%  '
Set-TextValue 'D3' '1.888.79'
Set-TextValue 'E3' '  -1.65%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '0.7734'
Set-TextValue 'E5' '  -4.24%  '
Set-TextValue 'D6' '244.65'
Set-TextValue 'E6' '  +0.03%  '
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'D8' '0.3139'
Set-TextValue 'E8' '  -3.95%  '
Set-TextValue 'D9' '0.07453'
Set-TextValue 'E9' '  +2.36%  '
Set-TextValue 'D10' '25.30'
Set-TextValue 'E10' '  -6.38%  '
Set-TextValue 'D11' '0.08123'
Set-TextValue 'E11' '  +0.44%  '
Set-TextValue 'D12' '0.7652'
Set-TextValue 'E12' '  -2.93%  '
Set-TextValue 'D13' '5.456'
Set-TextValue 'E13' '  +0.87%  '
Set-TextValue 'D14' '1.902.45'
Set-TextValue 'E14' '  -1.39%  '
Set-TextValue 'D15' '92.19'
Set-TextValue 'E15' '  -1.92%  '
Set-TextValue 'D16' '6.167'
Set-TextValue 'E16' '  +1.35%  '
Set-TextValue 'D17' '29.868.60'
Set-TextValue 'E17' '  -1.29%  '
Set-TextValue 'D18' '13.96'
Set-TextValue 'E18' '  -2.10%  '
Set-TextValue 'D19' '244.13'
Set-TextValue 'E19' '  -2.47%  '
Set-TextValue 'D20' '0.000007871'
Set-TextValue 'E20' '  -0.06%  '
Set-TextValue 'E21' '  +0.06%  '
Set-TextValue 'D22' '8.090'
Set-TextValue 'E22' '  -1.92%  '
Set-TextValue 'D23' '2.126.35'
Set-TextValue 'E23' '  -1.78%  '
Set-TextValue 'D24' '1.001'
Set-TextValue 'E24' '  -0.05%  '
Set-TextValue 'D25' '0.1572'
Set-TextValue 'E25' '  -4.30%  '
Set-TextValue 'D26' '9.419'
Set-TextValue 'E26' '  -0.74%  '
Set-TextValue 'D27' '162.69'
Set-TextValue 'E27' '  -2.98%  '
Set-TextValue 'D28' '18.79'
Set-TextValue 'E28' '  -1.26%  '
Set-TextValue 'D29' '2.041'
Set-TextValue 'E29' '  -5.56%  '
Set-TextValue 'D30' '1.438'
Set-TextValue 'E30' '  +3.24%  '
Set-TextValue 'D31' '1.549'
Set-TextValue 'E31' '  -0.27%  '
Set-TextValue 'E32' '  +1.64%  '
Set-TextValue 'D33' '4.092'
Set-TextValue 'E33' '  -1.59%  '
Set-TextValue 'D34' '0.05505'
Set-TextValue 'E34' '  -3.47%  '
Set-TextValue 'D35' '1.245'
Set-TextValue 'E35' '  -4.21%  '
Set-TextValue 'D36' '0.7584'
Set-TextValue 'E36' '  +0.88%  '
Set-TextValue 'D37' '1.003'
Set-TextValue 'E37' '  -0.11%  '
Set-TextValue 'D38' '2.647'
Set-TextValue 'E38' '  -2.96%  '
Set-TextValue 'D39' '0.01925'
Set-TextValue 'E39' '  -1.77%  '
Set-TextValue 'D40' '2.790'
Set-TextValue 'E40' '  -1.27%  '
Set-TextValue 'D41' '1.164.26'
Set-TextValue 'E41' '  +11.97%  '
Set-TextValue 'D42' '0.4451'
Set-TextValue 'E42' '  -1.95%  '
Set-TextValue 'D43' '73.85'
Set-TextValue 'E43' '  -0.61%  '
Set-TextValue 'D44' '5.963'
Set-TextValue 'E44' '  -1.04%  '
Set-TextValue 'D45' '0.8467'
Set-TextValue 'E45' '  -1.39%  '
Set-TextValue 'E46' '  +0.02%  '
Set-TextValue 'E47' '  -1.88%  '
Set-TextValue 'D48' '102.23'
Set-TextValue 'E48' '  -1.35%  '
Set-TextValue 'D49' '9.909'
Set-TextValue 'E49' '  -1.25%  '
Set-TextValue 'D50' '3.092'
Set-TextValue 'E50' '  -0.89%  '
Set-TextValue 'D51' '7.541'
Set-TextValue 'E51' '  -1.29%  '
